$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the bold/bordered/centered header style: reset header row (A1:E1) to the default style.
$ws.Range("A1:E1").ClearFormats()

# Add new row 5 of data
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 3
